$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.04184805085077817
$ws.Range("C2").Value = 0.00182260482758006

$ws.Range("B3").Value = 0.1413147618737045
$ws.Range("C3").Value = 0.001846019295771839

$ws.Range("B4").Value = 0.1225837217679183
$ws.Range("C4").Value = 0.002975274088728154

$ws.Range("B5").Value = 0.04293327166706364
$ws.Range("C5").Value = 0.001357542417309844

$ws.Range("B6").Value = 0.03119915676325567
$ws.Range("C6").Value = 0.001273838233403459
